# Add final data set for workflow 7: fill in candidate run times for rows 327-351 (E:J)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 327: E327=149259 F327=159806 G327=154333 H327=154281 I327=151909 J327=158719
$ws.Cells.Item(327, 5).Value = 149259
$ws.Cells.Item(327, 6).Value = 159806
$ws.Cells.Item(327, 7).Value = 154333
$ws.Cells.Item(327, 8).Value = 154281
$ws.Cells.Item(327, 9).Value = 151909
$ws.Cells.Item(327, 10).Value = 158719

# Row 328: E328=160216 F328=159071 G328=159982 H328=151660 I328=155168 J328=161745
$ws.Cells.Item(328, 5).Value = 160216
$ws.Cells.Item(328, 6).Value = 159071
$ws.Cells.Item(328, 7).Value = 159982
$ws.Cells.Item(328, 8).Value = 151660
$ws.Cells.Item(328, 9).Value = 155168
$ws.Cells.Item(328, 10).Value = 161745

# Row 329: E329=155096 F329=155892 G329=159195 H329=161682 I329=151608 J329=155097
$ws.Cells.Item(329, 5).Value = 155096
$ws.Cells.Item(329, 6).Value = 155892
$ws.Cells.Item(329, 7).Value = 159195
$ws.Cells.Item(329, 8).Value = 161682
$ws.Cells.Item(329, 9).Value = 151608
$ws.Cells.Item(329, 10).Value = 155097

# Row 330: E330=151339 F330=150637 G330=151507 H330=159123 I330=156141 J330=155171
$ws.Cells.Item(330, 5).Value = 151339
$ws.Cells.Item(330, 6).Value = 150637
$ws.Cells.Item(330, 7).Value = 151507
$ws.Cells.Item(330, 8).Value = 159123
$ws.Cells.Item(330, 9).Value = 156141
$ws.Cells.Item(330, 10).Value = 155171

# Row 331: E331=155159 F331=154328 G331=158971 H331=149918 I331=158024 J331=149561
$ws.Cells.Item(331, 5).Value = 155159
$ws.Cells.Item(331, 6).Value = 154328
$ws.Cells.Item(331, 7).Value = 158971
$ws.Cells.Item(331, 8).Value = 149918
$ws.Cells.Item(331, 9).Value = 158024
$ws.Cells.Item(331, 10).Value = 149561

# Row 332: E332=156713 F332=151664 G332=156271 H332=148155 I332=152090 J332=159538
$ws.Cells.Item(332, 5).Value = 156713
$ws.Cells.Item(332, 6).Value = 151664
$ws.Cells.Item(332, 7).Value = 156271
$ws.Cells.Item(332, 8).Value = 148155
$ws.Cells.Item(332, 9).Value = 152090
$ws.Cells.Item(332, 10).Value = 159538

# Row 333: E333=152646 F333=148195 G333=151259 H333=158289 I333=154263 J333=157163
$ws.Cells.Item(333, 5).Value = 152646
$ws.Cells.Item(333, 6).Value = 148195
$ws.Cells.Item(333, 7).Value = 151259
$ws.Cells.Item(333, 8).Value = 158289
$ws.Cells.Item(333, 9).Value = 154263
$ws.Cells.Item(333, 10).Value = 157163

# Row 334: E334=150416 F334=153587 G334=160497 H334=157690 I334=156956 J334=151569
$ws.Cells.Item(334, 5).Value = 150416
$ws.Cells.Item(334, 6).Value = 153587
$ws.Cells.Item(334, 7).Value = 160497
$ws.Cells.Item(334, 8).Value = 157690
$ws.Cells.Item(334, 9).Value = 156956
$ws.Cells.Item(334, 10).Value = 151569

# Row 335: E335=150558 F335=153484 G335=151469 H335=157312 I335=152320 J335=153058
$ws.Cells.Item(335, 5).Value = 150558
$ws.Cells.Item(335, 6).Value = 153484
$ws.Cells.Item(335, 7).Value = 151469
$ws.Cells.Item(335, 8).Value = 157312
$ws.Cells.Item(335, 9).Value = 152320
$ws.Cells.Item(335, 10).Value = 153058

# Row 336: E336=160636 F336=155126 G336=150385 H336=155922 I336=157596 J336=161473
$ws.Cells.Item(336, 5).Value = 160636
$ws.Cells.Item(336, 6).Value = 155126
$ws.Cells.Item(336, 7).Value = 150385
$ws.Cells.Item(336, 8).Value = 155922
$ws.Cells.Item(336, 9).Value = 157596
$ws.Cells.Item(336, 10).Value = 161473

# Row 337: E337=150645 F337=159451 G337=150554 H337=153991 I337=153751 J337=161354
$ws.Cells.Item(337, 5).Value = 150645
$ws.Cells.Item(337, 6).Value = 159451
$ws.Cells.Item(337, 7).Value = 150554
$ws.Cells.Item(337, 8).Value = 153991
$ws.Cells.Item(337, 9).Value = 153751
$ws.Cells.Item(337, 10).Value = 161354

# Row 338: E338=156562 F338=153898 G338=151209 H338=148625 I338=157186 J338=159051
$ws.Cells.Item(338, 5).Value = 156562
$ws.Cells.Item(338, 6).Value = 153898
$ws.Cells.Item(338, 7).Value = 151209
$ws.Cells.Item(338, 8).Value = 148625
$ws.Cells.Item(338, 9).Value = 157186
$ws.Cells.Item(338, 10).Value = 159051

# Row 339: E339=152697 F339=148624 G339=154348 H339=155957 I339=161553 J339=151272
$ws.Cells.Item(339, 5).Value = 152697
$ws.Cells.Item(339, 6).Value = 148624
$ws.Cells.Item(339, 7).Value = 154348
$ws.Cells.Item(339, 8).Value = 155957
$ws.Cells.Item(339, 9).Value = 161553
$ws.Cells.Item(339, 10).Value = 151272

# Row 340: E340=159805 F340=158494 G340=149645 H340=157835 I340=155737 J340=148135
$ws.Cells.Item(340, 5).Value = 159805
$ws.Cells.Item(340, 6).Value = 158494
$ws.Cells.Item(340, 7).Value = 149645
$ws.Cells.Item(340, 8).Value = 157835
$ws.Cells.Item(340, 9).Value = 155737
$ws.Cells.Item(340, 10).Value = 148135

# Row 341: E341=151137 F341=153214 G341=150485 H341=149302 I341=155232 J341=152814
$ws.Cells.Item(341, 5).Value = 151137
$ws.Cells.Item(341, 6).Value = 153214
$ws.Cells.Item(341, 7).Value = 150485
$ws.Cells.Item(341, 8).Value = 149302
$ws.Cells.Item(341, 9).Value = 155232
$ws.Cells.Item(341, 10).Value = 152814

# Row 342: E342=154038 F342=148593 G342=154376 H342=149363 I342=158208 J342=156615
$ws.Cells.Item(342, 5).Value = 154038
$ws.Cells.Item(342, 6).Value = 148593
$ws.Cells.Item(342, 7).Value = 154376
$ws.Cells.Item(342, 8).Value = 149363
$ws.Cells.Item(342, 9).Value = 158208
$ws.Cells.Item(342, 10).Value = 156615

# Row 343: E343=150410 F343=153397 G343=161361 H343=151654 I343=149936 J343=156155
$ws.Cells.Item(343, 5).Value = 150410
$ws.Cells.Item(343, 6).Value = 153397
$ws.Cells.Item(343, 7).Value = 161361
$ws.Cells.Item(343, 8).Value = 151654
$ws.Cells.Item(343, 9).Value = 149936
$ws.Cells.Item(343, 10).Value = 156155

# Row 344: E344=154457 F344=159251 G344=154559 H344=157733 I344=156178 J344=153571
$ws.Cells.Item(344, 5).Value = 154457
$ws.Cells.Item(344, 6).Value = 159251
$ws.Cells.Item(344, 7).Value = 154559
$ws.Cells.Item(344, 8).Value = 157733
$ws.Cells.Item(344, 9).Value = 156178
$ws.Cells.Item(344, 10).Value = 153571

# Row 345: E345=148663 F345=160185 G345=157611 H345=158705 I345=156744 J345=156570
$ws.Cells.Item(345, 5).Value = 148663
$ws.Cells.Item(345, 6).Value = 160185
$ws.Cells.Item(345, 7).Value = 157611
$ws.Cells.Item(345, 8).Value = 158705
$ws.Cells.Item(345, 9).Value = 156744
$ws.Cells.Item(345, 10).Value = 156570

# Row 346: E346=161155 F346=160133 G346=152867 H346=152467 I346=159315 J346=160696
$ws.Cells.Item(346, 5).Value = 161155
$ws.Cells.Item(346, 6).Value = 160133
$ws.Cells.Item(346, 7).Value = 152867
$ws.Cells.Item(346, 8).Value = 152467
$ws.Cells.Item(346, 9).Value = 159315
$ws.Cells.Item(346, 10).Value = 160696

# Row 347: E347=156139 F347=159442 G347=149280 H347=158300 I347=156599 J347=150548
$ws.Cells.Item(347, 5).Value = 156139
$ws.Cells.Item(347, 6).Value = 159442
$ws.Cells.Item(347, 7).Value = 149280
$ws.Cells.Item(347, 8).Value = 158300
$ws.Cells.Item(347, 9).Value = 156599
$ws.Cells.Item(347, 10).Value = 150548

# Row 348: E348=161378 F348=149368 G348=155626 H348=149459 I348=160088 J348=151771
$ws.Cells.Item(348, 5).Value = 161378
$ws.Cells.Item(348, 6).Value = 149368
$ws.Cells.Item(348, 7).Value = 155626
$ws.Cells.Item(348, 8).Value = 149459
$ws.Cells.Item(348, 9).Value = 160088
$ws.Cells.Item(348, 10).Value = 151771

# Row 349: E349=152869 F349=151286 G349=150760 H349=154107 I349=153650 J349=149388
$ws.Cells.Item(349, 5).Value = 152869
$ws.Cells.Item(349, 6).Value = 151286
$ws.Cells.Item(349, 7).Value = 150760
$ws.Cells.Item(349, 8).Value = 154107
$ws.Cells.Item(349, 9).Value = 153650
$ws.Cells.Item(349, 10).Value = 149388

# Row 350: E350=148842 F350=151933 G350=151485 H350=155333 I350=150901 J350=154840
$ws.Cells.Item(350, 5).Value = 148842
$ws.Cells.Item(350, 6).Value = 151933
$ws.Cells.Item(350, 7).Value = 151485
$ws.Cells.Item(350, 8).Value = 155333
$ws.Cells.Item(350, 9).Value = 150901
$ws.Cells.Item(350, 10).Value = 154840

# Row 351: E351=159202 F351=158282 G351=151242 H351=150645 I351=152393 J351=156971
$ws.Cells.Item(351, 5).Value = 159202
$ws.Cells.Item(351, 6).Value = 158282
$ws.Cells.Item(351, 7).Value = 151242
$ws.Cells.Item(351, 8).Value = 150645
$ws.Cells.Item(351, 9).Value = 152393
$ws.Cells.Item(351, 10).Value = 156971

# Update the active selection to match the author's final cursor position
$ws.Range("J33").Select()
